# Added account block test case
# Updates the "Test Results" sheet: refreshes execution timestamps for the
# "User Login with Invalid Credentials" steps, inserts a new step row for
# "Account Blocked After Multiple Wrong OTPs", and appends the trailing
# "Verify Home Page Loads Successfully" row that gets pushed down as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - re-run timestamp for the valid-login test case
$ws.Range("E2").Value = "28/03/2025 03:09:19 PM"

# Row 3 - Invalid Email Attempt
$ws.Range("B3").Value = ""
$ws.Range("E3").Value = "28/03/2025 03:09:22 PM"

# Row 4 - Not Registered Email Attempt
$ws.Range("B4").Value = ""
$ws.Range("E4").Value = "28/03/2025 03:09:23 PM"

# Row 5 - Invalid OTP Attempt
$ws.Range("B5").Value = ""
$ws.Range("E5").Value = "28/03/2025 03:09:25 PM"

# Row 6 - now the new "Account Blocked After Multiple Wrong OTPs" step
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "Account Blocked After Multiple Wrong OTPs"
$ws.Range("E6").Value = "28/03/2025 03:09:30 PM"
$ws.Range("F6").Value = "You have reached the maximum login attempts for the day. Please try again after 24 hours."

# Row 7 - "Click Go to Sign In" step shifts down from row 6
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "Click Go to Sign In"
$ws.Range("E7").Value = "28/03/2025 03:09:30 PM"
$ws.Range("F7").Value = "Navigated back to Get OTP page successfully"

# Row 8 - "No Sub-Steps" summary row for the invalid-credentials test case
$ws.Range("B8").Value = "User Login with Invalid Credentials"
$ws.Range("E8").Value = "28/03/2025 03:09:30 PM"

# Row 9 (new) - "Verify Home Page Loads Successfully" test case, pushed down
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Verify Home Page Loads Successfully"
$ws.Range("C9").Value = "No Sub-Steps"
$ws.Range("D9").Value = "PASSED"
$ws.Range("E9").Value = "28/03/2025 03:09:30 PM"
$ws.Range("F9").Value = "Test executed successfully."
